# Jueves dia grande, grabación coche, control de cambios
#
# Fill in the hours worked for 2023-04-26 (row 81) and 2023-04-27 (row 82)
# for each of the four people tracked in the "Hoja1" timesheet (columns
# G=Óscar, H=Inés, I=David, J=Daniel). The dependent totals/salary
# formulas further down the sheet (rows 139 and 142) recompute
# automatically from these new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 81 -> 2023-04-26 (serial 45042)
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 4
$ws.Range("I81").Value = 2
$ws.Range("J81").Value = 2

# Row 82 -> 2023-04-27 (serial 45043)
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 4
$ws.Range("I82").Value = 4
$ws.Range("J82").Value = 4

# Leave the sheet with the same cell focused as the author's last edit.
$ws.Range("K82").Select() | Out-Null
